$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / naturally non-numeric string updates
$ws.Range('D2').Value = '41.740.91'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '2.203.35'
$ws.Range('E3').Value = '  -0.29%  '
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('E5').Value = '  -1.71%  '
$ws.Range('E6').Value = '  -1.11%  '
$ws.Range('E7').Value = '  -2.26%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +5.83%  '
$ws.Range('E10').Value = '  +0.94%  '
$ws.Range('E11').Value = '  +2.50%  '
$ws.Range('E12').Value = '  -2.43%  '
$ws.Range('E13').Value = '  -2.81%  '
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('D15').Value = '2.529.76'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('E16').Value = '  -2.36%  '
$ws.Range('E17').Value = '  -2.61%  '
$ws.Range('D18').Value = '2.191.10'
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('D19').Value = '41.611.92'
$ws.Range('E19').Value = '  +0.31%  '
$ws.Range('D20').Value = '0.0₃0956'
$ws.Range('E20').Value = '  -0.61%  '
$ws.Range('E21').Value = '  -1.76%  '
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('E23').Value = '  -1.18%  '
$ws.Range('E24').Value = '  -1.38%  '
$ws.Range('E25').Value = '  -1.39%  '
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('E27').Value = '  -4.56%  '
$ws.Range('E28').Value = '  -5.08%  '
$ws.Range('E29').Value = '  -2.13%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('E30').Value = '  -2.37%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E31').Value = '  -2.45%  '
$ws.Range('E32').Value = '  -1.89%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('E33').Value = '  +7.08%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('E34').Value = '  -1.64%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('E35').Value = '  +4.54%  '
$ws.Range('E36').Value = '  -1.23%  '
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('E37').Value = '  -1.95%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E38').Value = '  +2.70%  '
$ws.Range('E39').Value = '  -1.93%  '
$ws.Range('E40').Value = '  +2.34%  '
$ws.Range('E41').Value = '  -1.35%  '
$ws.Range('E42').Value = '  +0.90%  '
$ws.Range('B43').Value = 'THORChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('E43').Value = '  -4.13%  '
$ws.Range('B44').Value = 'FTXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('E44').Value = '  +4.03%  '
$ws.Range('E45').Value = '  -3.97%  '
$ws.Range('E46').Value = '  -3.61%  '
$ws.Range('E47').Value = '  -2.62%  '
$ws.Range('E48').Value = '  -1.80%  '
$ws.Range('E49').Value = '  -0.28%  '
$ws.Range('E50').Value = '  -1.39%  '
$ws.Range('E51').Value = '  +6.36%  '

# Numeric-looking strings that must remain text (force text format, then clear style residue)
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').ClearFormats()
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.81'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.621'
$ws.Range('D6').ClearFormats()
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '68.08'
$ws.Range('D7').ClearFormats()
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.627'
$ws.Range('D9').ClearFormats()
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.56'
$ws.Range('D10').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '59.83'
$ws.Range('D11').ClearFormats()
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0936'
$ws.Range('D12').ClearFormats()
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.04'
$ws.Range('D13').ClearFormats()
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.104'
$ws.Range('D14').ClearFormats()
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.56'
$ws.Range('D16').ClearFormats()
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.866'
$ws.Range('D17').ClearFormats()
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.17'
$ws.Range('D21').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.22'
$ws.Range('D22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '231.40'
$ws.Range('D23').ClearFormats()
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.05'
$ws.Range('D24').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.90'
$ws.Range('D25').ClearFormats()
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.34'
$ws.Range('D27').ClearFormats()
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.42'
$ws.Range('D28').ClearFormats()
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.69'
$ws.Range('D29').ClearFormats()
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.14'
$ws.Range('D30').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '167.15'
$ws.Range('D31').ClearFormats()
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.34'
$ws.Range('D32').ClearFormats()
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0792'
$ws.Range('D33').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.121'
$ws.Range('D34').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.86'
$ws.Range('D35').ClearFormats()
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '26.10'
$ws.Range('D37').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.16'
$ws.Range('D38').ClearFormats()
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.60'
$ws.Range('D39').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0309'
$ws.Range('D40').ClearFormats()
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.22'
$ws.Range('D41').ClearFormats()
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '12.25'
$ws.Range('D42').ClearFormats()
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.63'
$ws.Range('D43').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.14'
$ws.Range('D44').ClearFormats()
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.89'
$ws.Range('D45').ClearFormats()
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.197'
$ws.Range('D46').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.56'
$ws.Range('D47').ClearFormats()
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0999'
$ws.Range('D48').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.16'
$ws.Range('D50').ClearFormats()
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.88'
$ws.Range('D51').ClearFormats()
